$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ideal-format")

# B3 used to be a hard-coded input (17041997); it is now derived from the
# breakdown rows 16-22 (same pattern already used elsewhere in the sheet).
$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"

# C16 used to be computed off a flat 66,000,000 total; it now reads the
# actual C3 total so it stays consistent if C3 changes.
$ws.Range("C16").Formula = "=C3-C17-C18-C19-C20-C21-C22"

# Selection moved from D11 to B1 in the saved view state.
$ws.Range("B1").Select()

$wb.Save()
